$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that currently sits at the very end
#    of the document (right after "活动最终结束权归记分组所有。"). Word
#    keeps this bookmark hidden from Bookmarks.Count / enumeration, but it
#    is still reachable (and deletable) by name. Do this *before* adding the
#    replacement bookmark below so the name is unambiguous (only one
#    "_GoBack" exists in the document at this point).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Insert a brand-new, completely empty paragraph at the very start of the
#    body that contains only an (empty) "_GoBack" bookmark. We build it via
#    InsertXML on a collapsed range at the start of the document so that no
#    stray run/formatting gets cloned onto the new paragraph (which is what
#    happens if InsertParagraphBefore() is used instead).
# ---------------------------------------------------------------------------
$goBackParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$startRange = $d.Range(0, 0)
$startRange.InsertXML($goBackParaXml)

# ---------------------------------------------------------------------------
# 3) Drop the w:proofErr spellStart/spellEnd pair that wraps the FIRST "git"
#    occurrence (the one right after "大家可以看到，在评分细则中，"). The two
#    later "git" occurrences keep their proofErr markers untouched. We
#    replace the whole containing paragraph with a corrected copy via
#    InsertXML so only that one proofErr pair disappears.
# ---------------------------------------------------------------------------
$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Execute("git占了相当大的比例") | Out-Null
$gitParagraph = $target.Paragraphs(1)

$fixedParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00C4520A" w:rsidRDefault="00C4520A" w:rsidP="006D4D53"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>大家可以看到，在评分细则中，</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>git</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>占了相当大的比例，这是我们最希望大家可以在这次</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>XLP</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>中学到的工具之一。大家一定要好好的阅读</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>的评分细则。其中我们比较关注的一点是——所有数据的源作者。希望你们可以对自己的作者身份负责，一方面是我们给分的重要凭据</w:t></w:r><w:r w:rsidR="00601101"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>之一，一方面这也是数字世界的凭据之一。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$gitParagraph.Range.InsertXML($fixedParaXml)
